$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AH, shifting old AH (Travel To) data to AI
$ws.Columns("AH").Insert()

# Set the new column header
$ws.Range("AH1").Value = "External File Links"

# Populate the wormhole.app external file links for the relevant rows
$ws.Range("AH2").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH3").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH4").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH5").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH6").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH7").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH8").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH9").Value = "https://wormhole.app/ROd2rN#xS3mraya3_qNbIe5hC6UTw"
$ws.Range("AH14").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
$ws.Range("AH19").Value = "https://wormhole.app/6Y1Z4E#TaAarl-rE6bznRDXWf3Okg"
